# Scheduled-runner style refresh of market-price derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) on several
# per-job "Leve Profits" sheets, plus removal of stale price rows whose
# items could no longer be priced on WVR.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# ALC
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value  = 2507.3333
$ws.Range("I94").Value  = 2400.7144
$ws.Range("K94").Value  = 2400.7144
$ws.Range("M94").Value  = -1949.7144

$ws.Range("H131").Value = 1012.8
$ws.Range("I131").Value = 891.44446
$ws.Range("J131").Value = 2105
$ws.Range("K131").Value = 2674.33338
$ws.Range("L131").Value = 6315
$ws.Range("M131").Value = 2365.66662
$ws.Range("N131").Value = -16395

# ---------------------------------------------------------------
# BSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1773.5869
$ws.Range("I105").Value = 1610.7407
$ws.Range("J105").Value = 2005
$ws.Range("K105").Value = 1610.7407
$ws.Range("L105").Value = 2005
$ws.Range("M105").Value = 136.2592999999999
$ws.Range("N105").Value = -5499

# ---------------------------------------------------------------
# CRP
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value  = 946.0714
$ws.Range("I16").Value  = 603.25
$ws.Range("J16").Value  = 1403.1666
$ws.Range("K16").Value  = 603.25
$ws.Range("L16").Value  = 1403.1666
$ws.Range("M16").Value  = -316.25
$ws.Range("N16").Value  = -1977.1666

$ws.Range("H113").Value = 946.0714
$ws.Range("I113").Value = 603.25
$ws.Range("J113").Value = 1403.1666
$ws.Range("K113").Value = 603.25
$ws.Range("L113").Value = 1403.1666
$ws.Range("M113").Value = 1566.75
$ws.Range("N113").Value = -5743.1666

$ws.Range("H141").Value = 54608.668
$ws.Range("J141").Value = 54608.668
$ws.Range("L141").Value = 54608.668
$ws.Range("N141").Value = -64968.668

# ---------------------------------------------------------------
# CUL
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value  = 1761.5222
$ws.Range("I68").Value  = 1490.7675
$ws.Range("J68").Value  = 2009.234
$ws.Range("K68").Value  = 4472.3025
$ws.Range("L68").Value  = 6027.701999999999
$ws.Range("M68").Value  = -3661.3025
$ws.Range("N68").Value  = -7649.701999999999

$ws.Range("H71").Value  = 1761.5222
$ws.Range("I71").Value  = 1490.7675
$ws.Range("J71").Value  = 2009.234
$ws.Range("K71").Value  = 13416.9075
$ws.Range("L71").Value  = 18083.106
$ws.Range("M71").Value  = -9360.907499999999
$ws.Range("N71").Value  = -26195.106

# ---------------------------------------------------------------
# GSM
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 541.53125
$ws.Range("I107").Value = 504.30435
$ws.Range("J107").Value = 636.6667
$ws.Range("K107").Value = 504.30435
$ws.Range("L107").Value = 636.6667
$ws.Range("M107").Value = 1415.69565
$ws.Range("N107").Value = -4476.6667

$ws.Range("H132").Value = 72552.44500000001
$ws.Range("I132").Value = 81520.88
$ws.Range("J132").Value = 16499.75
$ws.Range("K132").Value = 244562.64
$ws.Range("L132").Value = 49499.25
$ws.Range("M132").Value = -242032.64
$ws.Range("N132").Value = -54559.25

# ---------------------------------------------------------------
# LTW
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6250.383
$ws.Range("I132").Value = 8368.585999999999
$ws.Range("J132").Value = 2837.7222
$ws.Range("K132").Value = 25105.758
$ws.Range("L132").Value = 8513.1666
$ws.Range("M132").Value = -22575.758
$ws.Range("N132").Value = -13573.1666

# ---------------------------------------------------------------
# WVR — drop the now-unpriceable market columns (H:N) for the rows
# whose items fell out of the Universalis price cache this run.
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
$wvrRows = 119,120,121,122,123,124,125,126,127,128,129,130,131,132,133,135,136,137,138,139,140,141
foreach ($r in $wvrRows) {
    $ws.Range("H" + $r + ":N" + $r).ClearContents()
}
